# Usuarios workbook update:
#  - row 3's phone-like columns (A, C) switch from text to plain numbers.
#  - a new row 4 is added for a new user (Wilson Fabio), whose phone-like
#    columns (A, C) are entered as text, matching the original row 3 layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: A3/C3 become genuine numbers instead of text -------------------
$ws.Range("A3").Value = 102020120
$ws.Range("C3").Value = 30343323222

# --- Row 4: brand new user record -------------------------------------------
# A4 keeps the same bold/bordered/centered style as A2/A3 (style index 1).
# The leading "'" forces Excel to store the numeric-looking value as text
# (shared string) instead of silently converting it to a number.
$ws.Range("A4").Value = "'1020810810"
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats -> re-apply A3's style (1)

$ws.Range("B4").Value = "Wilson Fabio"

# C4 has no special style (same as C2/C3), so after forcing it to text we
# paste the formatting of an untouched, default-styled cell back onto it.
$ws.Range("C4").Value = "'31200029299"
$ws.Range("Z1").Copy()
$ws.Range("C4").PasteSpecial(-4122)  # xlPasteFormats -> back to default style (0)

$ws.Range("D4").Value = "w@f.com"
